# Applies the "prepared demo without implementation of panel brightness and
# sleep after" edit:
#   - Typography sheet: add a new "Fifteen" typography row (row 8) mirroring
#     the existing Default/Large/Small/Medium rows.
#   - Translation sheet: fill in the three previously-blank rows (33-35)
#     with new demo text entries (software version, IP address, visual
#     style) instead of leaving them empty.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Typography sheet: row 8
# ---------------------------------------------------------------------
$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("B8").Value = "Fifteen"
$wsTypo.Range("C8").Value = "verdana.ttf"
$wsTypo.Range("D8").Value = 15
$wsTypo.Range("E8").Value = 4
$wsTypo.Range("F8").Value = "?"
$wsTypo.Range("G8").Value = "abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ0123456789"
$wsTypo.Range("H8").Value = "0-9,A-Z,a-z"

# Columns B:E carry a non-default column style, but the existing data rows
# (4-7) use the workbook default ("Normal") style on each cell - reset the
# new cells to match so row 8 looks like the others.
$wsTypo.Range("B8:E8").Style = "Normal"

# Column I is left blank for every typography row, but still has an actual
# (empty) cell record - touch it so row 8 keeps the same shape as rows 4-7.
$wsTypo.Range("I8").Style = "Normal"

# ---------------------------------------------------------------------
# Translation sheet: rows 33-35
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B33").Value = "SingleUseId39"
$wsTrans.Range("C33").Value = "Medium"
$wsTrans.Range("D33").Value = "Left"
$wsTrans.Range("E33").Value = "Software v0.5.1"
$wsTrans.Range("F33").Value = "LTR"

$wsTrans.Range("B34").Value = "SingleUseId40"
$wsTrans.Range("C34").Value = "Medium"
$wsTrans.Range("D34").Value = "Right"
$wsTrans.Range("E34").Value = "IP: 192.168.9.24"
$wsTrans.Range("F34").Value = "LTR"

$wsTrans.Range("B35").Value = "SingleUseId41"
$wsTrans.Range("C35").Value = "Default"
$wsTrans.Range("D35").Value = "Left"
$wsTrans.Range("E35").Value = "Visual Style:"
$wsTrans.Range("F35").Value = "LTR"

# Columns B:F on this sheet also carry a non-default column style; reset
# the newly written cells to the workbook default, matching the other
# populated rows in the table.
$wsTrans.Range("B33:F35").Style = "Normal"
